# Applies the cryptos.xlsx price/volume update described in the commit diff.
# Values are written with a leading apostrophe to force Excel to store them as
# literal text (matching the source inlineStr cells) instead of auto-converting
# number-looking strings (e.g. "0.730", "54.70") into numeric values that would
# drop the trailing zero. The style is reset to "Normal" right after so the
# cell doesn't pick up a lingering quote-prefix style.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($row, $col, $value) {
    $cell = $ws.Cells.Item($row, $col)
    $cell.Value = "'" + $value
    $cell.Style = "Normal"
}

Set-TextValue 2 4 "65.891.23"
Set-TextValue 2 5 "  +1.48%  "
Set-TextValue 3 4 "3.174.17"
Set-TextValue 3 5 "  +0.71%  "
Set-TextValue 4 5 "  +0.16%  "
Set-TextValue 5 4 "593.46"
Set-TextValue 5 5 "  +3.50%  "
Set-TextValue 6 4 "152.34"
Set-TextValue 6 5 "  +2.00%  "
Set-TextValue 7 5 "  +0.06%  "
Set-TextValue 8 4 "3.176.34"
Set-TextValue 8 5 "  +0.76%  "
Set-TextValue 9 4 "0.533"
Set-TextValue 9 5 "  +1.50%  "
Set-TextValue 10 4 "0.158"
Set-TextValue 10 5 "  -0.66%  "
Set-TextValue 11 4 "6.06"
Set-TextValue 11 5 "  -0.63%  "
Set-TextValue 12 4 "0.509"
Set-TextValue 12 5 "  +2.64%  "
Set-TextValue 13 4 "0.0000266"
Set-TextValue 13 5 "  +1.02%  "
Set-TextValue 14 4 "38.52"
Set-TextValue 14 5 "  +4.20%  "
Set-TextValue 15 4 "3.707.19"
Set-TextValue 15 5 "  +1.13%  "
Set-TextValue 16 4 "66.007.24"
Set-TextValue 16 5 "  +1.58%  "
Set-TextValue 17 4 "7.36"
Set-TextValue 17 5 "  +3.90%  "
Set-TextValue 18 4 "3.188.20"
Set-TextValue 18 5 "  +0.58%  "
Set-TextValue 19 5 "  +0.29%  "
Set-TextValue 20 4 "506.19"
Set-TextValue 20 5 "  +0.12%  "
Set-TextValue 21 4 "15.19"
Set-TextValue 21 5 "  +2.70%  "
Set-TextValue 22 4 "0.730"
Set-TextValue 22 5 "  +1.92%  "
Set-TextValue 23 4 "7.94"
Set-TextValue 23 5 "  +3.20%  "
Set-TextValue 24 4 "14.86"
Set-TextValue 24 5 "  -2.41%  "
Set-TextValue 25 4 "84.35"
Set-TextValue 25 5 "  +0.29%  "
Set-TextValue 26 5 "  -0.08%  "
Set-TextValue 27 4 "9.19"
Set-TextValue 27 5 "  +4.51%  "
Set-TextValue 28 4 "2.97"
Set-TextValue 28 5 "  +2.45%  "
Set-TextValue 29 4 "2.27"
Set-TextValue 29 5 "  +5.22%  "
Set-TextValue 30 4 "6.89"
Set-TextValue 30 5 "  +11.91%  "
Set-TextValue 31 4 "2.85"
Set-TextValue 31 5 "  +2.45%  "
Set-TextValue 32 4 "28.02"
Set-TextValue 32 5 "  +1.70%  "
Set-TextValue 33 5 "  +2.92%  "
Set-TextValue 34 5 "  +0.29%  "
Set-TextValue 35 4 "6.46"
Set-TextValue 35 5 "  -0.30%  "
Set-TextValue 36 4 "54.70"
Set-TextValue 36 5 "  -0.25%  "
Set-TextValue 37 4 "0.0891"
Set-TextValue 37 5 "  -0.54%  "
Set-TextValue 38 4 "478.63"
Set-TextValue 38 5 "  +3.34%  "
Set-TextValue 39 4 "0.0414"
Set-TextValue 39 5 "  -1.16%  "
Set-TextValue 40 4 "8.76"
Set-TextValue 40 5 "  +1.43%  "
Set-TextValue 41 4 "2.85"
Set-TextValue 41 5 "  -4.51%  "
Set-TextValue 42 4 "0.121"
Set-TextValue 42 5 "  +3.63%  "
Set-TextValue 43 4 "0.295"
Set-TextValue 43 5 "  +4.52%  "
Set-TextValue 44 2 "Maker"
Set-TextValue 44 3 "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
Set-TextValue 44 4 "2.901.65"
Set-TextValue 44 5 "  -4.75%  "
Set-TextValue 45 2 "PEPE"
Set-TextValue 45 3 "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
Set-TextValue 45 4 "0.0₃0637"
Set-TextValue 45 5 "  +9.38%  "
Set-TextValue 46 4 "2.38"
Set-TextValue 46 5 "  -2.03%  "
Set-TextValue 47 4 "28.11"
Set-TextValue 47 5 "  -1.25%  "
Set-TextValue 48 5 "  -0.03%  "
Set-TextValue 49 5 "  +1.45%  "
Set-TextValue 50 4 "2.28"
Set-TextValue 50 5 "  +2.20%  "
Set-TextValue 51 4 "2.57"
Set-TextValue 51 5 "  +5.11%  "
